# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (matching commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1)
$ws1.Range("F2").Value = 1334
$ws1.Range("F4").Value = 14629
$ws1.Range("F5").Value = 17809
$ws1.Range("F6").Value = 147
$ws1.Range("F7").Value = 74
$ws1.Range("F10").Value = 29
$ws1.Range("F16").Value = 54
$ws1.Range("F17").Value = 157
$ws1.Range("F19").Value = 1340
$ws1.Range("F24").Value = 7355
$ws1.Range("F25").Value = 981
$ws1.Range("F26").Value = 7
$ws1.Range("F27").Value = 41
$ws1.Range("F28").Value = 1181
$ws1.Range("F30").Value = 5875
$ws1.Range("F33").Value = 144
$ws1.Range("F35").Value = 231
$ws1.Range("F36").Value = 5118

# Sheet "全部类型" (sheet4)
$ws4.Range("F2").Value = 1334
$ws4.Range("F4").Value = 14629
$ws4.Range("F5").Value = 17809
$ws4.Range("F6").Value = 147
$ws4.Range("F7").Value = 74
$ws4.Range("F10").Value = 29
$ws4.Range("F16").Value = 54
$ws4.Range("F17").Value = 157
$ws4.Range("F19").Value = 1340
$ws4.Range("F25").Value = 7355
$ws4.Range("F26").Value = 981
$ws4.Range("F27").Value = 7
$ws4.Range("F28").Value = 41
$ws4.Range("F29").Value = 1181
$ws4.Range("F32").Value = 5875
$ws4.Range("F35").Value = 144
$ws4.Range("F37").Value = 231
$ws4.Range("F38").Value = 5118
